$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 30, shifting rows 30-127 down to 31-128
$ws.Rows("30:30").Insert()

# Populate the new row 30 by copying most fields from row 31 (which used to be
# row 30 before the shift), then overwrite the fields that actually changed.
for ($c = 1; $c -le 18; $c++) {
    $ws.Cells.Item(30, $c).Value = $ws.Cells.Item(31, $c).Value2
}

# Overwrite the fields that differ per the diff
$ws.Cells.Item(30, 4).Value = 45133   # D30 Fecha
$ws.Cells.Item(30, 10).Value = 100    # J30 Volumen
$ws.Cells.Item(30, 11).Value = 4500   # K30 Precio minimo
$ws.Cells.Item(30, 12).Value = 5000   # L30 Precio maximo
$ws.Cells.Item(30, 13).Value = 4750   # M30 Precio promedio ponderado
$ws.Cells.Item(30, 16).Value = 132    # P30 Precio $/Kg
